$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.894.66'
$ws.Range("E2").Value = '  +0.89%  '

# Row 3
$ws.Range("D3").Value = '2.212.96'
$ws.Range("E3").Value = '  +0.39%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.62%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.15'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.66%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.630'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.85'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.06'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0940'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.26%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("E14").Value = '  -0.35%  '

# Row 15
$ws.Range("D15").Value = '2.540.71'
$ws.Range("E15").Value = '  +0.22%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.71%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.873'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.17%  '

# Row 18
$ws.Range("D18").Value = '2.199.24'
$ws.Range("E18").Value = '  +0.91%  '

# Row 19
$ws.Range("D19").Value = '41.817.59'
$ws.Range("E19").Value = '  +0.76%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0960'
$ws.Range("E20").Value = '  +0.52%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.01%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.45'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.64%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '

# Row 28
$ws.Range("E28").Value = '  -5.56%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.95%  '

# Row 30
$ws.Range("E30").Value = '  -2.17%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.42'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.00%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0799'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.93%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.62%  '

# Row 35
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.121'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.38%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.123'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.62%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.41%  '

# Row 38
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.63'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.24%  '

# Row 39
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.79%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0310'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.97%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.12%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.84%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.75%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.200'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.46%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '

# Row 48
$ws.Range("E48").Value = '  -1.02%  '

# Row 49
$ws.Range("E49").Value = '  -0.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '

# Row 51
$ws.Range("E51").Value = '  +6.39%  '
